$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsGroups = $wb.Worksheets.Item("Groups")

# --- Login sheet (sheet1): update expected-result text, drop the now-unused
#     "Incorrect login or password" cell for the bad-login test row ---
$wsLogin.Range("D2").Value = "Admin Admin (admin)"
$wsLogin.Range("D3").Value = "LabManager LabManager (lab manager)"
$wsLogin.Range("D4").ClearContents()

# --- Groups sheet (sheet3): rename the test group used in the "addGroup" case ---
$wsGroups.Range("D4").Value = "OstrTestGroup5"

# --- Selection / active-tab bookkeeping ---
# Groups tab loses its previous selection/active state, Login tab becomes
# the active tab with D4 selected.
$wsGroups.Range("A32:A33").Select()
$wsLogin.Range("D4").Select()
